$wb = $excel.ActiveWorkbook

# --- Practitioners sheet ---
$wsP = $wb.Worksheets.Item("Practitioners")
$wsP.Activate() | Out-Null
$wsP.Columns.Item(1).ColumnWidth = 13.8
$wsP.Columns.Item(3).ColumnWidth = 12.2
$wsP.Columns.Item(6).ColumnWidth = 12

$wsP.Range("A6").Value = "PHN999:NFP02"
$wsP.Range("B6").Value = "P01"
$wsP.Range("C6").Value = 8
$wsP.Range("D6").Value = 1
$wsP.Range("E6").Value = 1973
$wsP.Range("F6").Value = 2
$wsP.Range("G6").Value = 1
$wsP.Range("H6").Value = 1
$wsP.Range("I6").Value = "tag1"

$wsP.Columns.Item(7).Select() | Out-Null

# --- Service Contacts sheet (left as the active tab, matching the source workbook) ---
$wsSC = $wb.Worksheets.Item("Service Contacts")
$wsSC.Activate() | Out-Null
$wsSC.Columns.Item(1).ColumnWidth = 13.666666
$wsSC.Range("D3").Select() | Out-Null
